$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 53
$ws.Range("H53").Value = 1368
$ws.Range("I53").Value = 240.14285
$ws.Range("K53").Value = 240.14285
$ws.Range("M53").Value = 396.85715
# Row 69
$ws.Range("H69").Value = 20785.643
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 20785.643
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 62356.929
$ws.Range("M69").ClearContents()
$ws.Range("N69").Value = -64104.929
# Row 72
$ws.Range("H72").Value = 20785.643
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 20785.643
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 187070.787
$ws.Range("M72").ClearContents()
$ws.Range("N72").Value = -195806.787
# Row 87
$ws.Range("H87").Value = 38272.2
$ws.Range("J87").Value = 39191.332
$ws.Range("L87").Value = 39191.332
$ws.Range("N87").Value = -41687.332
# Row 90
$ws.Range("H90").Value = 38272.2
$ws.Range("J90").Value = 39191.332
$ws.Range("L90").Value = 117573.996
$ws.Range("N90").Value = -130053.996
# Row 111
$ws.Range("H111").Value = 3713.7273
$ws.Range("I111").Value = 3356.5
$ws.Range("K111").Value = 10069.5
$ws.Range("M111").Value = -7002.5

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 1641.2222
$ws.Range("I32").Value = 1641.2222
$ws.Range("K32").Value = 1641.2222
$ws.Range("M32").Value = -1354.2222
# Row 74
$ws.Range("H74").Value = 24393606
$ws.Range("I74").Value = 29414710
$ws.Range("K74").Value = 29414710
$ws.Range("M74").Value = -29413836
# Row 77
$ws.Range("H77").Value = 24393606
$ws.Range("I77").Value = 29414710
$ws.Range("K77").Value = 147073550
$ws.Range("M77").Value = -147069182
# Row 122
$ws.Range("H122").Value = 7507.5835
$ws.Range("I122").Value = 5917.364
$ws.Range("K122").Value = 17752.092
$ws.Range("M122").Value = -15302.092

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 9
$ws.Range("H9").Value = 20051
$ws.Range("I9").Value = 20051
$ws.Range("K9").Value = 20051
$ws.Range("M9").Value = -19883
# Row 22
$ws.Range("H22").Value = 4765104.5
$ws.Range("J22").Value = 7938007
$ws.Range("L22").Value = 7938007
$ws.Range("N22").Value = -7938353
# Row 44
$ws.Range("H44").Value = 30022
$ws.Range("I44").Value = 20045
$ws.Range("J44").Value = 39999
$ws.Range("K44").Value = 20045
$ws.Range("L44").Value = 39999
$ws.Range("M44").Value = -19548
$ws.Range("N44").Value = -40993
# Row 63
$ws.Range("H63").Value = 49999
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 49999
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 49999
$ws.Range("M63").ClearContents()
$ws.Range("N63").Value = -51371
# Row 66
$ws.Range("H66").Value = 49999
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 49999
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 149997
$ws.Range("M66").ClearContents()
$ws.Range("N66").Value = -156861
# Row 123
$ws.Range("H123").Value = 89999
$ws.Range("J123").Value = 89999
$ws.Range("L123").Value = 89999
$ws.Range("N123").Value = -99799

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 4635.143
$ws.Range("I31").Value = 3147
$ws.Range("J31").Value = 6619.3335
$ws.Range("K31").Value = 3147
$ws.Range("L31").Value = 6619.3335
$ws.Range("M31").Value = -2852
$ws.Range("N31").Value = -7209.3335
# Row 34
$ws.Range("H34").Value = 4635.143
$ws.Range("I34").Value = 3147
$ws.Range("J34").Value = 6619.3335
$ws.Range("K34").Value = 3147
$ws.Range("L34").Value = 6619.3335
$ws.Range("M34").Value = -2945
$ws.Range("N34").Value = -7023.3335
# Row 57
$ws.Range("H57").Value = 0
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("M57").ClearContents()
$ws.Range("N57").ClearContents()
# Row 86
$ws.Range("H86").Value = 13924.8
$ws.Range("I86").Value = 10750
$ws.Range("K86").Value = 10750
$ws.Range("M86").Value = -9627
# Row 89
$ws.Range("H89").Value = 13924.8
$ws.Range("I89").Value = 10750
$ws.Range("K89").Value = 53750
$ws.Range("M89").Value = -48134
# Row 99
$ws.Range("H99").Value = 3246.5881
$ws.Range("I99").Value = 3131.3635
$ws.Range("J99").Value = 3457.8333
$ws.Range("K99").Value = 3131.3635
$ws.Range("L99").Value = 3457.8333
$ws.Range("M99").Value = -1633.3635
$ws.Range("N99").Value = -6453.8333
# Row 108
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()
# Row 122
$ws.Range("H122").Value = 3741.6924
$ws.Range("I122").Value = 3575.3333
$ws.Range("J122").Value = 4116
$ws.Range("K122").Value = 10725.9999
$ws.Range("L122").Value = 12348
$ws.Range("M122").Value = -8275.999899999999
$ws.Range("N122").Value = -17248
# Row 126
$ws.Range("H126").Value = 3246.5881
$ws.Range("I126").Value = 3131.3635
$ws.Range("J126").Value = 3457.8333
$ws.Range("K126").Value = 9394.0905
$ws.Range("L126").Value = 10373.4999
$ws.Range("M126").Value = -6924.0905
$ws.Range("N126").Value = -15313.4999

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 12
$ws.Range("H12").Value = 95.92856999999999
$ws.Range("I12").Value = 37.75
$ws.Range("K12").Value = 113.25
$ws.Range("M12").Value = 59.75
# Row 122
$ws.Range("H122").Value = 1500
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 3107.3572
$ws.Range("I80").Value = 3428.2222
$ws.Range("K80").Value = 3428.2222
$ws.Range("M80").Value = -2430.2222
# Row 83
$ws.Range("H83").Value = 3107.3572
$ws.Range("I83").Value = 3428.2222
$ws.Range("K83").Value = 17141.111
$ws.Range("M83").Value = -12149.111
# Row 122
$ws.Range("H122").Value = 5022.846
$ws.Range("I122").Value = 3441.4167
$ws.Range("J122").Value = 24000
$ws.Range("K122").Value = 10324.2501
$ws.Range("L122").Value = 72000
$ws.Range("M122").Value = -7874.250100000001
$ws.Range("N122").Value = -76900
# Row 126
$ws.Range("H126").Value = 3150.9524
$ws.Range("I126").Value = 3084.0715
$ws.Range("K126").Value = 9252.2145
$ws.Range("M126").Value = -6782.2145
# Row 134
$ws.Range("H134").Value = 61499
$ws.Range("J134").Value = 61499
$ws.Range("L134").Value = 184497
$ws.Range("N134").Value = -189567

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 2374.4211
$ws.Range("I22").Value = 2855.375
$ws.Range("J22").Value = 2024.6364
$ws.Range("K22").Value = 2855.375
$ws.Range("L22").Value = 2024.6364
$ws.Range("M22").Value = -2560.375
$ws.Range("N22").Value = -2614.6364
# Row 27
$ws.Range("H27").Value = 2374.4211
$ws.Range("I27").Value = 2855.375
$ws.Range("J27").Value = 2024.6364
$ws.Range("K27").Value = 2855.375
$ws.Range("L27").Value = 2024.6364
$ws.Range("M27").Value = -2748.375
$ws.Range("N27").Value = -2238.6364
# Row 40
$ws.Range("H40").Value = 4457.647
$ws.Range("I40").Value = 4291.385
$ws.Range("K40").Value = 4291.385
$ws.Range("M40").Value = -4155.385
# Row 46
$ws.Range("H46").Value = 3962.5
$ws.Range("I46").Value = 2740
$ws.Range("J46").Value = 6000
$ws.Range("K46").Value = 2740
$ws.Range("L46").Value = 6000
$ws.Range("M46").Value = -2552
$ws.Range("N46").Value = -6376
# Row 58
$ws.Range("H58").Value = 14658.286
$ws.Range("I58").Value = 8768.166999999999
$ws.Range("J58").Value = 49999
$ws.Range("K58").Value = 8768.166999999999
$ws.Range("L58").Value = 49999
$ws.Range("M58").Value = -8508.166999999999
$ws.Range("N58").Value = -50519
# Row 122
$ws.Range("H122").Value = 31599.4
$ws.Range("J122").Value = 50000
$ws.Range("L122").Value = 150000
$ws.Range("N122").Value = -154900

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 5
$ws.Range("H5").Value = 54750
$ws.Range("J5").Value = 54750
$ws.Range("L5").Value = 54750
$ws.Range("N5").Value = -54974
# Row 100
$ws.Range("H100").Value = 1758.6923
$ws.Range("I100").Value = 2415.875
$ws.Range("J100").Value = 707.2
$ws.Range("K100").Value = 4831.75
$ws.Range("L100").Value = 1414.4
$ws.Range("M100").Value = -4290.75
$ws.Range("N100").Value = -2496.4
